# ============================================================
# Add "2022-Q4" sheet with fund-holding data, insert a
# corresponding summary row on "总计", and renumber the
# existing quarter index values.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- 1. Update the "总计" (summary) sheet: sheet1 ----------------------
$total = $wb.Worksheets.Item(1)

# Insert a new row 2 (shifts 2022-Q3/Q2/Q1/2021-Q4 rows down by one)
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

# Give the new A2 the same style as the other index cells (bold/center/border)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 26
$total.Range("D2").Value = 2.47

# Renumber the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# ---- 2. Insert a brand-new "2022-Q4" worksheet right after "总计" -------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Pull an existing quarter sheet (now holding the former "2022-Q3" data)
# to copy its header/row formatting from, since COM doesn't expose raw
# style indices directly.
$styleSrc = $wb.Worksheets.Item("2022-Q3")

$fundData = @(
    ,(0, "506000", "南方科创板 3 年定开混合", "24.35", "91.59", "3.09", "0.7524", 10)
    ,(1, "481010", "工银中小盘混合", "15.71", "91.30", "2.92", "0.4587", 10)
    ,(2, "161610", "融通领先成长混合（LOF）A", "13.83", "88.35", "1.50", "0.2074", 10)
    ,(3, "011404", "融通鑫新成长混合C", "5.57", "94.42", "2.45", "0.1365", 9)
    ,(4, "161601", "融通新蓝筹混合", "11.21", "70.40", "1.19", "0.1334", 9)
    ,(5, "240017", "华宝新兴产业混合", "2.74", "86.64", "4.60", "0.1260", 3)
    ,(6, "010114", "华宝新兴成长混合A", "3.09", "89.74", "3.97", "0.1227", 4)
    ,(7, "008811", "鹏华科技创新混合", "2.83", "89.03", "3.53", "0.0999", 6)
    ,(8, "015553", "融通价值成长混合A", "3.54", "94.64", "2.80", "0.0991", 10)
    ,(9, "001192", "上投摩根整合驱动灵活配置混合A", "3.19", "92.87", "3.07", "0.0979", 3)
    ,(10, "011403", "融通鑫新成长混合A", "2.22", "94.42", "2.45", "0.0544", 9)
    ,(11, "519929", "长信电子信息行业量化灵活配置混合A", "0.84", "90.40", "4.90", "0.0412", 8)
    ,(12, "015554", "融通价值成长混合C", "1.36", "94.64", "2.80", "0.0381", 10)
    ,(13, "005593", "上投摩根创新商业模式灵活配置混合A", "0.94", "94.26", "4.04", "0.0380", 6)
    ,(14, "002281", "建信裕利灵活配置混合", "0.85", "93.43", "3.16", "0.0269", 10)
    ,(15, "002955", "融通新趋势灵活配置混合", "0.74", "90.61", "1.61", "0.0119", 9)
    ,(16, "005382", "泰康睿利量化多策略混合C", "0.43", "89.86", "1.98", "0.0085", 6)
    ,(17, "005381", "泰康睿利量化多策略混合A", "0.41", "89.86", "1.98", "0.0081", 6)
    ,(18, "010646", "融通价值趋势混合A", "0.41", "74.27", "1.30", "0.0053", 9)
    ,(19, "010647", "融通价值趋势混合C", "0.10", "74.27", "1.30", "0.0013", 9)
    ,(20, "001708", "东兴改革精选灵活配置混合", "0.03", "89.54", "2.50", "0.0008", 10)
    ,(21, "009241", "融通领先成长混合（LOF）C", "0.04", "88.35", "1.50", "0.0006", 10)
    ,(22, "017197", "华宝新兴成长混合C", "0.00", "89.74", "3.97", $null, 4)
    ,(23, "013153", "长信电子信息行业量化灵活配置混合C", "0.00", "90.40", "4.90", $null, 8)
    ,(24, "016418", "上投摩根创新商业模式灵活配置混合C", "0.00", "94.26", "4.04", $null, 6)
    ,(25, "016920", "上投摩根整合驱动灵活配置混合C", "0.00", "92.87", "3.07", $null, 3)
)

# ---- 3. Header row ------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).NumberFormat = "@"
    $q4.Cells.Item(1, $col).Value = $headers[$col - 2]
    $q4.Cells.Item(1, $col).ClearFormats()
}
$styleSrc.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# ---- 4. Data rows ---------------------------------------------------------
$r = 2
foreach ($row in $fundData) {
    $q4.Cells.Item($r, 1).Value = $row[0]

    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 2).ClearFormats()

    $q4.Cells.Item($r, 3).NumberFormat = "@"
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 3).ClearFormats()

    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 4).ClearFormats()

    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 5).ClearFormats()

    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 6).ClearFormats()

    if ($row[6] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).NumberFormat = "@"
        $q4.Cells.Item($r, 7).Value = $row[6]
        $q4.Cells.Item($r, 7).ClearFormats()
    }

    $q4.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# Column-A index cells (A2:A27) get the bold/center/border style, same as
# the other quarter sheets.
$styleSrc.Range("A2").Copy()
$q4.Range("A2:A27").PasteSpecial(-4122)

$q4.Range("A1").Select()
